# Insert a new price-record row for Murcott / Primera (row 372) into the
# "Vega Modelo de Temuco - Mandarina" sheet. This pushes the existing rows
# 372-468 down to 373-469 (dimension grows from A1:T468 to A1:T469).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 372; everything below shifts down by one.
$ws.Rows(372).Insert()

# Populate the new row with the new price record.
$ws.Range("A372").Value = 10
$ws.Range("B372").Value = "Vega Modelo de Temuco"
$ws.Range("C372").Value = "La Araucanía"
$ws.Range("D372").Value = 44642
$ws.Range("E372").Value = 9
$ws.Range("F372").Value = "Fruta"
$ws.Range("G372").Value = 100102
$ws.Range("H372").Value = "Cítricos"
$ws.Range("I372").Value = 100102004
$ws.Range("J372").Value = "Mandarina"
$ws.Range("K372").Value = "Murcott"
$ws.Range("L372").Value = "Primera"
$ws.Range("M372").Value = 80
$ws.Range("N372").Value = 19800
$ws.Range("O372").Value = 19800
$ws.Range("P372").Value = 19800
$ws.Range("Q372").Value = "$/bandeja 18 kilos"
$ws.Range("R372").Value = "Región de O'Higgins"
$ws.Range("S372").Value = 1100
$ws.Range("T372").Value = 18
